$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45992
$ws.Range("B2").Value = 87.02
$ws.Range("C2").Value = 75.78
$ws.Range("D2").Value = 70.5
$ws.Range("E2").Value = 70.16
$ws.Range("F2").Value = 77.31
$ws.Range("G2").Value = 88.73999999999999
$ws.Range("H2").Value = 99.67
$ws.Range("I2").Value = 108.35
$ws.Range("J2").Value = 113.51
$ws.Range("K2").Value = 95.20999999999999
$ws.Range("L2").Value = 75.12
$ws.Range("M2").Value = 68.2
$ws.Range("N2").Value = 70.8
$ws.Range("O2").Value = 69.56999999999999
$ws.Range("P2").Value = 69.88
$ws.Range("Q2").Value = 73.38
$ws.Range("R2").Value = 85.65000000000001
$ws.Range("S2").Value = 95.12
$ws.Range("T2").Value = 104.87
$ws.Range("U2").Value = 122.96
$ws.Range("V2").Value = 137.87
$ws.Range("W2").Value = 138.25
$ws.Range("X2").Value = 109.01
$ws.Range("Y2").Value = 95.78
$ws.Range("Z2").Value = 91.78
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 120.23
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 138.06
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 113.92
$ws.Range("AG2").Value = "0h-16h"
